$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# "Ready for handoff" -> "In Translation" (status text changed on the
# Overview sheet as well as each per-language handoff sheet)
$ws1.Range("E2").Value = "In Translation"
$ws1.Range("F2").Value = "In Translation"
$ws2.Range("C2").Value = "In Translation"
$ws3.Range("C2").Value = "In Translation"

# The Status columns narrow to fit the shorter replacement text.
$ws1.Columns.Item(5).ColumnWidth = 13.4101845877511
$ws1.Columns.Item(6).ColumnWidth = 13.4101845877511
$ws2.Columns.Item(3).ColumnWidth = 13.4101845877511
$ws3.Columns.Item(3).ColumnWidth = 13.4101845877511
